$d = $word.ActiveDocument

# Change 1: "outbreak occurring.  " -> "outbreak on the horizon  "
$d.Content.Find.Execute("warming areas which show signs of an outbreak occurring.  ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "warming areas which show signs of an outbreak on the horizon  ", 2)

# Change 2: "combination can lead to society getting back to " -> "combination can perhaps lower the risk profile that leads to society returning to "
$d.Content.Find.Execute("combination can lead to society getting back to ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "combination can perhaps lower the risk profile that leads to society returning to ", 2)
